$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the selection moving to E8 after the edit (matches author's saved state)
$ws.Range("E8").Select()
